$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source values (bill_length_mm, flipper_length_mm, species probabilities,
# model_version, prediction_timestamp) are all stored as text in this workbook,
# even the numeric-looking ones. Force the target range to Text format before
# writing so values land as shared strings instead of being auto-coerced to
# numbers, then drop back to the Normal style so no extra cell formatting is
# left behind on the written cells.
$ws.Range("A2:H3").NumberFormat = "@"

# Row 2: update the existing prediction with refreshed model output
$ws.Range("A2").Value = "39.1"
$ws.Range("B2").Value = "181.0"
$ws.Range("C2").Value = "Adelie"
$ws.Range("D2").Value = "1.0"
$ws.Range("E2").Value = "0.0"
$ws.Range("F2").Value = "0.0"
$ws.Range("G2").Value = "v1.0"
$ws.Range("H2").Value = "2025-05-04 20:48:45"

# Row 3: newly appended prediction row
$ws.Range("A3").Value = "46.5"
$ws.Range("B3").Value = "192.0"
$ws.Range("C3").Value = "Chinstrap"
$ws.Range("D3").Value = "0.09"
$ws.Range("E3").Value = "0.91"
$ws.Range("F3").Value = "0.0"
$ws.Range("G3").Value = "v1.0"
$ws.Range("H3").Value = "2025-05-04 20:48:45"

# Restore the default (Normal) cell style now that the text values are in place
$ws.Range("A2:H3").Style = "Normal"
